$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.365.96"
$ws.Range("E2").Value = "  +1.66%  "

$ws.Range("D3").Value = "1.833.45"
$ws.Range("E3").Value = "  +1.04%  "

$ws.Range("D4").Value = "'1.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.87%  "

$ws.Range("D5").Value = "'314.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.78%  "

$ws.Range("D6").Value = "'1.009"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.75%  "

$ws.Range("D7").Value = "'0.4742"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.98%  "

$ws.Range("E8").Value = "  +0.97%  "

$ws.Range("E9").Value = "  +1.37%  "

$ws.Range("D10").Value = "'0.8868"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.99%  "

$ws.Range("E11").Value = "  +1.00%  "

$ws.Range("D12").Value = "1.908.53"
$ws.Range("E12").Value = "  +7.04%  "

$ws.Range("D13").Value = "'0.07333"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.26%  "

$ws.Range("D14").Value = "'5.454"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").Value = "'93.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.95%  "

$ws.Range("D16").Value = "'6.587"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.35%  "

$ws.Range("D17").Value = "'1.009"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.55%  "

$ws.Range("D18").Value = "'0.000008820"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.36%  "

$ws.Range("D19").Value = "'1.009"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.79%  "

$ws.Range("D20").Value = "27.567.86"
$ws.Range("E20").Value = "  +2.40%  "

$ws.Range("D21").Value = "'14.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.19%  "

$ws.Range("E22").Value = "  +0.42%  "

$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").Value = "2.096.46"
$ws.Range("E24").Value = "  +2.09%  "

$ws.Range("D25").Value = "'1.895"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").Value = "'151.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.66%  "

$ws.Range("E27").Value = "  +1.64%  "

$ws.Range("D28").Value = "'2.143"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.06%  "

$ws.Range("D29").Value = "'5.249"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("D30").Value = "'117.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.69%  "

$ws.Range("D31").Value = "'0.09000"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.07%  "

$ws.Range("D32").Value = "'0.7553"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("D33").Value = "'1.179"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.62%  "

$ws.Range("D34").Value = "'4.551"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.61%  "

$ws.Range("D35").Value = "'2.944"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.44%  "

$ws.Range("E36").Value = "  +0.90%  "

$ws.Range("D37").Value = "'1.104"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.78%  "

$ws.Range("D38").Value = "'0.05355"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.36%  "

$ws.Range("D39").Value = "'0.01953"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.24%  "

$ws.Range("E40").Value = "  +0.54%  "

$ws.Range("D41").Value = "'7.319"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.00%  "

$ws.Range("D42").Value = "'2.398"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.75%  "

$ws.Range("D43").Value = "'0.5325"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.28%  "

$ws.Range("E44").Value = "  +0.50%  "

$ws.Range("D45").Value = "'8.496"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.83%  "

$ws.Range("D46").Value = "'0.4915"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").Value = "'10.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.78%  "

$ws.Range("D48").Value = "'105.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.91%  "

$ws.Range("E49").Value = "  +0.91%  "

$ws.Range("E50").Value = "  +0.96%  "

$ws.Range("D51").Value = "'0.06296"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.05%  "

